# Applies the cryptos list refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.638.20"
$ws.Range("E2").Value = "  -0.98%  "
$ws.Range("D3").Value = "2.563.54"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'578.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.96%  "
$ws.Range("D6").Value = "'143.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.94%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("E9").Value = "  -2.73%  "
$ws.Range("D10").Value = "'5.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.17%  "
$ws.Range("E11").Value = "  -0.69%  "
$ws.Range("D12").Value = "'0.349"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.43%  "
$ws.Range("D13").Value = "'26.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.03%  "
$ws.Range("D14").Value = "3.022.53"
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("D15").Value = "62.539.59"
$ws.Range("E15").Value = "  -0.96%  "
$ws.Range("D16").Value = "'0.0000144"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.60%  "
$ws.Range("D17").Value = "2.568.84"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").Value = "'11.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.52%  "
$ws.Range("D19").Value = "'337.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.85%  "
$ws.Range("E20").Value = "  -2.17%  "
$ws.Range("D21").Value = "'6.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.37%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "'66.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.88%  "
$ws.Range("D24").Value = "'1.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.20%  "
$ws.Range("E25").Value = "  +1.23%  "
$ws.Range("E26").Value = "  -4.77%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("D28").Value = "'7.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.83%  "
$ws.Range("E29").Value = "  -4.58%  "
$ws.Range("E30").Value = "  -2.65%  "
$ws.Range("D31").Value = "'454.55"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.79%  "
$ws.Range("E32").Value = "  -4.89%  "
$ws.Range("D33").Value = "'176.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.84%  "
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("E36").Value = "  -3.67%  "
$ws.Range("D37").Value = "'18.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.12%  "
$ws.Range("D38").Value = "'4.42"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.37%  "
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("E40").Value = "  -4.71%  "
$ws.Range("D41").Value = "'40.45"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.53%  "
$ws.Range("D42").Value = "'156.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.08%  "
$ws.Range("E43").Value = "  -4.61%  "
$ws.Range("D45").Value = "'20.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.36%  "
$ws.Range("E46").Value = "  -3.90%  "
$ws.Range("D47").Value = "'0.0955"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.39%  "
$ws.Range("E48").Value = "  -4.18%  "
$ws.Range("D49").Value = "'17.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.46%  "
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("E51").Value = "  -5.20%  "
